$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New header cells (row 1) + new data cells tied to them (row 2).
# The write order below matches the order new entries must land in the
# shared-strings table: identifier, GT152, taxnum, phonenum, faxnum.
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "identifier"
$ws.Range("E2").Value = "GT152"
$ws.Range("F1").Value = "taxnum"
$ws.Range("G1").Value = "phonenum"
$ws.Range("H1").Value = "faxnum"

$ws.Range("F2").Value = 15422
$ws.Range("G2").Value = 112455
$ws.Range("H2").Value = 1548754

# ---------------------------------------------------------------------
# The "Hyperlink" style previously sat on the (then-empty) G2/H2 cells.
# Those columns now hold real data, so clear their formatting back to
# Normal and move the blank styled placeholder cell over to I2 instead.
# ---------------------------------------------------------------------
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Style = "Normal"
$ws.Range("I2").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Column width bestFit adjustments for the new/resized columns.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 9.5
$ws.Columns.Item(5).ColumnWidth = 8.666666666666666
$ws.Columns.Item(6).ColumnWidth = 11.5
$ws.Columns.Item(7).ColumnWidth = 8.166666666666666
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666
$ws.Columns.Item(9).ColumnWidth = 22.166666666666668

# ---------------------------------------------------------------------
# Selection moves from A2 to H2.
# ---------------------------------------------------------------------
[void]$ws.Range("H2").Select()

# ---------------------------------------------------------------------
# The saved window is slightly shorter than before.
# ---------------------------------------------------------------------
[void]($excel.ActiveWindow.Height = 2175)
